# "Updated performance for V4.08.XX"
# Fill in the previously-blank performance metrics (Albedo/Shading/Shadow/RGB
# PSNR+SSIM pairs, columns B:I) for model versions V4.08.5 .. V4.08.8, which
# live on the "Maps 2 RGB" sheet in rows 33-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maps 2 RGB")

$data = @{
    33 = @(5.86, 0.60440000000000005, 12.8194, 0.79800000000000004, 14.534000000000001, 0.78539999999999999, 17.504999999999999, 0.81569999999999998)
    34 = @(7.0288000000000004, 0.6361, 12.803000000000001, 0.79959999999999998, 14.8584, 0.80179999999999996, 17.394500000000001, 0.83679999999999999)
    35 = @(6.2824999999999998, 0.61470000000000002, 13.659800000000001, 0.82269999999999999, 13.672800000000001, 0.77710000000000001, 17.3598, 0.81569999999999998)
    36 = @(6.4057000000000004, 0.63, 13.8142, 0.8196, 14.1731, 0.77769999999999995, 17.474599999999999, 0.81830000000000003)
}

foreach ($row in ($data.Keys | Sort-Object)) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # Column B is 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# Scroll the frozen pane up a bit and move the selection, matching the
# author's new cursor position after entering the data.
$window = $excel.ActiveWindow
$window.ScrollRow = 26
$ws.Range("I32").Select()
